$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, shifting existing rows 25-72 down to 26-73.
$ws.Rows("25:25").Insert()

# Populate the new row 25 with the new weekly data point.
$ws.Cells.Item(25, 1).Value = 4
$ws.Cells.Item(25, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(25, 3).Value = "Los Lagos"
$ws.Cells.Item(25, 4).Value = 45002
$ws.Cells.Item(25, 5).Value = 10
$ws.Cells.Item(25, 6).Value = 100112030
$ws.Cells.Item(25, 7).Value = "Poroto granado"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 60
$ws.Cells.Item(25, 11).Value = 38000
$ws.Cells.Item(25, 12).Value = 38000
$ws.Cells.Item(25, 13).Value = 38000
$ws.Cells.Item(25, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value = "Región del Maule"
$ws.Cells.Item(25, 16).Value = 1520
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"
